$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the yellow "needs review" highlight from the release-date cells
#    in column C for rows 28-34 and 36-39 (row 35 is already un-highlighted,
#    and its format - numFmtId 168 / fontId 1 / borderId 0, fillId 0 - is
#    exactly what the target cells should end up with). Copy its format
#    (format-only paste) onto the target cells so the existing "no fill"
#    style gets reused instead of minting a new one.
# ---------------------------------------------------------------------------
$ws.Range("C35").Copy()
$ws.Range("C28:C34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C35").Copy()
$ws.Range("C36:C39").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Roll the daily/weekly market-data history one column to the right
#    (R<-Q, S<-R, T<-S, U<-T) and drop in the new "as of" value in Q, along
#    with the refreshed "as of" date in N, for each updated series.
# ---------------------------------------------------------------------------

# Row 29 - 5yr,5yr Forward (T5YIFR): as-of date 2026-02-23 -> 2026-02-24
$oldQ = $ws.Range("Q29").Value2
$oldR = $ws.Range("R29").Value2
$oldS = $ws.Range("S29").Value2
$oldT = $ws.Range("T29").Value2
$ws.Range("N29").Value = 46077
$ws.Range("U29").Value = $oldT
$ws.Range("T29").Value = $oldS
$ws.Range("S29").Value = $oldR
$ws.Range("R29").Value = $oldQ
$ws.Range("Q29").Value = 2.12

# Row 30 - 10yr TIPS (T10YIE): as-of date 2026-02-23 -> 2026-02-24
$oldQ = $ws.Range("Q30").Value2
$oldR = $ws.Range("R30").Value2
$oldS = $ws.Range("S30").Value2
$oldT = $ws.Range("T30").Value2
$ws.Range("N30").Value = 46077
$ws.Range("U30").Value = $oldT
$ws.Range("T30").Value = $oldS
$ws.Range("S30").Value = $oldR
$ws.Range("R30").Value = $oldQ
$ws.Range("Q30").Value = 2.26

# Row 47 - FFR (DFF): as-of date 2026-02-20 -> 2026-02-23 (values unchanged, all 3.64)
$ws.Range("N47").Value = 46076

# Row 48 - 2y UST (DGS2): as-of date 2026-02-20 -> 2026-02-23
$oldQ = $ws.Range("Q48").Value2
$oldR = $ws.Range("R48").Value2
$oldS = $ws.Range("S48").Value2
$oldT = $ws.Range("T48").Value2
$ws.Range("N48").Value = 46076
$ws.Range("U48").Value = $oldT
$ws.Range("T48").Value = $oldS
$ws.Range("S48").Value = $oldR
$ws.Range("R48").Value = $oldQ
$ws.Range("Q48").Value = 3.43

# Row 49 - 5y UST (DGS5): as-of date 2026-02-20 -> 2026-02-23
$oldQ = $ws.Range("Q49").Value2
$oldR = $ws.Range("R49").Value2
$oldS = $ws.Range("S49").Value2
$oldT = $ws.Range("T49").Value2
$ws.Range("N49").Value = 46076
$ws.Range("U49").Value = $oldT
$ws.Range("T49").Value = $oldS
$ws.Range("S49").Value = $oldR
$ws.Range("R49").Value = $oldQ
$ws.Range("Q49").Value = 3.59

# Row 50 - 10y UST (DGS10): as-of date 2026-02-20 -> 2026-02-23
$oldQ = $ws.Range("Q50").Value2
$oldR = $ws.Range("R50").Value2
$oldS = $ws.Range("S50").Value2
$oldT = $ws.Range("T50").Value2
$ws.Range("N50").Value = 46076
$ws.Range("U50").Value = $oldT
$ws.Range("T50").Value = $oldS
$ws.Range("S50").Value = $oldR
$ws.Range("R50").Value = $oldQ
$ws.Range("Q50").Value = 4.03

# Row 52 - BAA (DBAA): as-of date 2026-02-20 -> 2026-02-23
$oldQ = $ws.Range("Q52").Value2
$oldR = $ws.Range("R52").Value2
$oldS = $ws.Range("S52").Value2
$oldT = $ws.Range("T52").Value2
$ws.Range("N52").Value = 46076
$ws.Range("U52").Value = $oldT
$ws.Range("T52").Value = $oldS
$ws.Range("S52").Value = $oldR
$ws.Range("R52").Value = $oldQ
$ws.Range("Q52").Value = 5.76
